$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 9737.526
$ws.Range("I21").Value = 5017
$ws.Range("J21").Value = 9999.777
$ws.Range("K21").Value = 5017
$ws.Range("L21").Value = 9999.777
$ws.Range("M21").Value = -4549
$ws.Range("N21").Value = -10935.777

$ws.Range("H23").Value = 9737.526
$ws.Range("I23").Value = 5017
$ws.Range("J23").Value = 9999.777
$ws.Range("K23").Value = 5017
$ws.Range("L23").Value = 9999.777
$ws.Range("M23").Value = -4783
$ws.Range("N23").Value = -10467.777

$ws.Range("H98").Value = 1156.1904
$ws.Range("I98").Value = 714
$ws.Range("J98").Value = 10000
$ws.Range("K98").Value = 714
$ws.Range("L98").Value = 10000
$ws.Range("M98").Value = 784
$ws.Range("N98").Value = -12996

$ws.Range("H122").Value = 1156.1904
$ws.Range("I122").Value = 714
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 2142
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = 308
$ws.Range("N122").Value = -34900

$ws.Range("H135").Value = 21586.877
$ws.Range("I135").Value = 28797.223
$ws.Range("J135").Value = 1619.7693
$ws.Range("K135").Value = 259175.007
$ws.Range("L135").Value = 14577.9237
$ws.Range("M135").Value = -256640.007
$ws.Range("N135").Value = -19647.9237

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3093.238
$ws.Range("I45").Value = 1304.6364
$ws.Range("J45").Value = 5060.7
$ws.Range("K45").Value = 1304.6364
$ws.Range("L45").Value = 5060.7
$ws.Range("M45").Value = -927.6364000000001
$ws.Range("N45").Value = -5814.7

$ws.Range("H61").Value = 1111.4678
$ws.Range("I61").Value = 1120.34
$ws.Range("J61").Value = 1074.5
$ws.Range("K61").Value = 1120.34
$ws.Range("L61").Value = 1074.5
$ws.Range("M61").Value = -908.3399999999999
$ws.Range("N61").Value = -1498.5

$ws.Range("H74").Value = 933.5
$ws.Range("I74").Value = 921.6
$ws.Range("J74").Value = 1052.5
$ws.Range("K74").Value = 921.6
$ws.Range("L74").Value = 1052.5
$ws.Range("M74").Value = -47.60000000000002
$ws.Range("N74").Value = -2800.5

$ws.Range("H77").Value = 933.5
$ws.Range("I77").Value = 921.6
$ws.Range("J77").Value = 1052.5
$ws.Range("K77").Value = 4608
$ws.Range("L77").Value = 5262.5
$ws.Range("M77").Value = -240
$ws.Range("N77").Value = -13998.5

$ws.Range("H132").Value = 96078.17
$ws.Range("I132").Value = 99130.25999999999
$ws.Range("J132").Value = 18250
$ws.Range("K132").Value = 297390.78
$ws.Range("L132").Value = 54750
$ws.Range("M132").Value = -294860.78
$ws.Range("N132").Value = -59810

$ws.Range("H136").Value = 1111.4678
$ws.Range("I136").Value = 1120.34
$ws.Range("J136").Value = 1074.5
$ws.Range("K136").Value = 3361.02
$ws.Range("L136").Value = 3223.5
$ws.Range("M136").Value = -811.0199999999995
$ws.Range("N136").Value = -8323.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H31").Value = 1383.8
$ws.Range("I31").Value = 1044.44
$ws.Range("J31").Value = 2232.2
$ws.Range("K31").Value = 1044.44
$ws.Range("L31").Value = 2232.2
$ws.Range("M31").Value = -749.4400000000001
$ws.Range("N31").Value = -2822.2

$ws.Range("H34").Value = 1383.8
$ws.Range("I34").Value = 1044.44
$ws.Range("J34").Value = 2232.2
$ws.Range("K34").Value = 1044.44
$ws.Range("L34").Value = 2232.2
$ws.Range("M34").Value = -842.4400000000001
$ws.Range("N34").Value = -2636.2

$ws.Range("H138").Value = 46107
$ws.Range("J138").Value = 46107
$ws.Range("L138").Value = 46107
$ws.Range("N138").Value = -56387

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 868.7143
$ws.Range("I92").Value = 329.66666
$ws.Range("J92").Value = 1273
$ws.Range("K92").Value = 988.9999799999999
$ws.Range("L92").Value = 3819
$ws.Range("M92").Value = 259.0000200000001
$ws.Range("N92").Value = -6315

$ws.Range("H107").Value = 519.5806
$ws.Range("I107").Value = 230.95653
$ws.Range("J107").Value = 1349.375
$ws.Range("K107").Value = 692.86959
$ws.Range("L107").Value = 4048.125
$ws.Range("M107").Value = 1227.13041
$ws.Range("N107").Value = -7888.125

$ws.Range("H136").Value = 3731.9285
$ws.Range("I136").Value = 1661.7059
$ws.Range("J136").Value = 4395.9624
$ws.Range("K136").Value = 4985.1177
$ws.Range("L136").Value = 13187.8872
$ws.Range("M136").Value = 114.8823000000002
$ws.Range("N136").Value = -23387.8872

$ws.Range("H138").Value = 2293.25
$ws.Range("I138").Value = 1220.2858
$ws.Range("J138").Value = 3127.7778
$ws.Range("K138").Value = 3660.8574
$ws.Range("L138").Value = 9383.3334
$ws.Range("M138").Value = 1479.1426
$ws.Range("N138").Value = -19663.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1607.6936
$ws.Range("I132").Value = 1310.6222
$ws.Range("J132").Value = 2394.0588
$ws.Range("K132").Value = 3931.8666
$ws.Range("L132").Value = 7182.176399999999
$ws.Range("M132").Value = -1401.8666
$ws.Range("N132").Value = -12242.1764

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1618
$ws.Range("I68").Value = 1600
$ws.Range("J68").Value = 1690
$ws.Range("K68").Value = 1600
$ws.Range("L68").Value = 1690
$ws.Range("M68").Value = -851
$ws.Range("N68").Value = -3188

$ws.Range("H71").Value = 1618
$ws.Range("I71").Value = 1600
$ws.Range("J71").Value = 1690
$ws.Range("K71").Value = 8000
$ws.Range("L71").Value = 8450
$ws.Range("M71").Value = -4256
$ws.Range("N71").Value = -15938

$ws.Range("H132").Value = 1501.871
$ws.Range("I132").Value = 1434.1694
$ws.Range("J132").Value = 2833.3333
$ws.Range("K132").Value = 4302.5082
$ws.Range("L132").Value = 8499.999899999999
$ws.Range("M132").Value = -1772.5082
$ws.Range("N132").Value = -13559.9999

$ws.Range("H136").Value = 1130.8591
$ws.Range("I136").Value = 978.0161000000001
$ws.Range("J136").Value = 2183.7778
$ws.Range("K136").Value = 2934.0483
$ws.Range("L136").Value = 6551.3334
$ws.Range("M136").Value = -384.0483000000004
$ws.Range("N136").Value = -11651.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 32754.75
$ws.Range("I24").Value = 12009
$ws.Range("J24").Value = 39670
$ws.Range("K24").Value = 12009
$ws.Range("L24").Value = 39670
$ws.Range("M24").Value = -11779
$ws.Range("N24").Value = -40130

$ws.Range("H136").Value = 1361.2031
$ws.Range("I136").Value = 1342.4906
$ws.Range("J136").Value = 1451.3636
$ws.Range("K136").Value = 4027.4718
$ws.Range("L136").Value = 4354.0908
$ws.Range("M136").Value = -1477.4718
$ws.Range("N136").Value = -9454.0908
